$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point is inserted at row 79 (same market/variety profile
# as the former row 79), pushing all subsequent rows down by one.
$ws.Rows(79).Insert()

$ws.Range("A79").Value = 11
$ws.Range("B79").Value = "Vega Monumental Concepción"
$ws.Range("C79").Value = "Bíobío"
$ws.Range("D79").Value = 44874
$ws.Range("E79").Value = 8
$ws.Range("F79").Value = 100112032
$ws.Range("G79").Value = "Zapallo italiano"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 220
$ws.Range("K79").Value = 7000
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = 7273
$ws.Range("N79").Value = "$/caja 50 unidades"
$ws.Range("O79").Value = "Región de O'Higgins"
$ws.Range("P79").Value = 145
$ws.Range("Q79").Value = 50
$ws.Range("R79").Value = "Hortaliza"
